# Update the time-slot labels in column C and move the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "22:40-22:45"
$ws.Range("C7").Value = "22:45-22:50"

$ws.Range("C11").Select()
